# "Decision final con TP y SL"
# Restructures the "Insumos" sheet (renames / regroups some inputs, adds a
# Bid-Ask Spread input, clears the old TP/SL hard/soft-limit rows) and adds
# a new "Zona_PG_Cierre" worksheet describing the Stop Loss / Low Loss /
# Low Profit / Take Profit decision zones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Rebuild "Insumos" from a clean slate so row numbers / merges exactly
#    match the new layout.
# ---------------------------------------------------------------------
$ws.Cells.UnMerge()
$ws.Cells.Clear()

# --- Section: Informacion para filtro de datos (row 1) ---
$ws.Range("A1:B1").Merge()
$ws.Range("A1").Value = "Informacion para filtro de datos"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Size = 12
$ws.Range("A1:B1").Borders.Item(9).LineStyle = 1
$ws.Range("A1:B1").HorizontalAlignment = -4108

$ws.Range("A2").Value = "N_FH_por_dia"
$ws.Range("B2").Value = 13
$ws.Range("A3").Value = "N_FH_Cierre_Descartadas"
$ws.Range("B3").Value = 3

# --- Section: Informacion para escalamiento de posición (row 5) ---
$ws.Range("A5:B5").Merge()
$ws.Range("A5").Value = "Informacion para escalamiento de posición"
$ws.Range("A5:B5").Font.Bold = $true
$ws.Range("A5:B5").Font.Size = 12
$ws.Range("A5:B5").Borders.Item(9).LineStyle = 1
$ws.Range("A5:B5").HorizontalAlignment = -4108

$ws.Range("A6").Value = "Apalancamiento máximo"
$ws.Range("B6").Value = 7
$ws.Range("B6").NumberFormat = "0.0"

# --- Section: Informacion de negociación (row 9) ---
$ws.Range("A9:B9").Merge()
$ws.Range("A9").Value = "Informacion de negociación"
$ws.Range("A9:B9").Font.Bold = $true
$ws.Range("A9:B9").Font.Size = 12
$ws.Range("A9:B9").Borders.Item(9).LineStyle = 1
$ws.Range("A9:B9").HorizontalAlignment = -4108

$ws.Range("A10").Value = "Bid-Ask Spread (BAS)"
$ws.Range("B10").Value = 0.0001
$ws.Range("B10").NumberFormat = "0.00%"
$ws.Range("B10").Font.Size = 11

$ws.Range("A11").Value = "Comisión"
$ws.Range("B11").Value = 0.05
$ws.Range("B11").NumberFormat = "0.0%"
$ws.Range("B11").Font.Size = 11

$ws.Range("B12:B14").NumberFormat = "0.0%"
$ws.Range("B12:B14").Font.Size = 11

# --- Section: Informacion para cálculo de variables de riesgo (row 16) ---
$ws.Range("A16:B16").Merge()
$ws.Range("A16").Value = "Informacion para cálculo de variables de riesgo"
$ws.Range("A16:B16").Font.Bold = $true
$ws.Range("A16:B16").Font.Size = 12
$ws.Range("A16:B16").Borders.Item(9).LineStyle = 1
$ws.Range("A16:B16").HorizontalAlignment = -4108

$ws.Range("A17").Value = "VentanaMovilVol (Dias)"
$ws.Range("B17").Value = 90

$ws.Range("A18").Value = "Significancia"
$ws.Range("B18").Value = 0.05
$ws.Range("B18").NumberFormat = "0.0%"
$ws.Range("B18").Font.Size = 11

# --- Section: Informacion para mostrar en gráficos (row 20) ---
$ws.Range("A20:B20").Merge()
$ws.Range("A20").Value = "Informacion para mostrar en gráficos"
$ws.Range("A20:B20").Font.Bold = $true
$ws.Range("A20:B20").Font.Size = 12
$ws.Range("A20:B20").Borders.Item(9).LineStyle = 1
$ws.Range("A20:B20").HorizontalAlignment = -4108

$ws.Range("A21").Value = "U_MDD_Objetivo"
$ws.Range("B21").Value = 3
$ws.Range("B21").NumberFormat = "0.0"

$ws.Range("A22").Value = "RA_MDD_Objetivo"
$ws.Range("B22").Value = 0.2

$ws.Range("A23").Value = "Sharpe_Objetivo"
$ws.Range("B23").Value = 0.2

$ws.Range("A24").Value = "Sortino_Objetivo"
$ws.Range("B24").Value = 0.4

# ---------------------------------------------------------------------
# 2. Add the "Zona_PG_Cierre" worksheet describing the TP / SL zones.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws)
$ws2.Name = "Zona_PG_Cierre"

$ws2.Range("A1").Value = "Zona_PG"
$ws2.Range("B1").Value = "MinimoRazon"
$ws2.Range("C1").Value = "MaximoRazon"
$ws2.Range("D1").Value = "Cierre"
$ws2.Range("A1:D1").Font.Bold = $true
$ws2.Range("A1:D1").Borders.Item(9).LineStyle = 1
$ws2.Range("A1:D1").HorizontalAlignment = -4108

$ws2.Range("A2").Value = "Stop Loss"
$ws2.Range("C2").Value = -0.025
$ws2.Range("D2").Value = 1

$ws2.Range("A3").Value = "Low Loss"
$ws2.Range("B3").Value = -0.025
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 0.5

$ws2.Range("A4").Value = "Low Profit"
$ws2.Range("B4").Value = 0
$ws2.Range("C4").Value = 0.025
$ws2.Range("D4").Value = 0

$ws2.Range("A5").Value = "Take Profit"
$ws2.Range("B5").Value = 0.025
$ws2.Range("D5").Value = 0.5

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B3:C5").NumberFormat = "0.0%"
$ws2.Range("C2").NumberFormat = "0.0%"
$ws2.Range("D2:D5").NumberFormat = "0%"
$ws2.Range("B2:D5").HorizontalAlignment = -4108

$ws2.Range("A1").EntireColumn.ColumnWidth = 9.88671875
$ws2.Range("B1:C1").EntireColumn.ColumnWidth = 12.6640625

$ws.Select()
